# Fix bug in back polygon config file
#
# The "kind" condition values in column A of the C_BackPolygon.conf sheet
# were stored as bare tokens (e.g. "kind=0136"); the fix quotes the value
# portion of every kind=... clause (e.g. "kind="0136"") so the config
# parser treats it as a quoted literal. Rows that use "id=..." (the
# Shanghai Disney / Guangzhou Changlong rows) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_BackPolygon.conf")

$lastRow = $ws.Range("A1").End(-4121).Row  # xlDown
if ($lastRow -lt 79) { $lastRow = 79 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("A$r")
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*kind=*") {
        $newval = $val -replace 'kind=([^;\s|()]+)', 'kind="$1"'
        if ($newval -ne $val) {
            $cell.Value = $newval
        }
    }
}

# Reflect the author's final selection/active cell on the first sheet.
$ws.Select()
$ws.Range("A71").Select()
